$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "Tasks Assigned" section: swap the two volunteer names between the
# "Delivery html, css, js" bullet and the "Restaurant html, css, js"
# bullet, and drop the now-empty spacer paragraph that used to sit
# right after the "Restaurant" bullet.
# ------------------------------------------------------------------

# --- Delivery html, css, js -> ... bullet --------------------------
# "Dea Hasanaj, Brend Zmijanej"  =>  "Danja Korreshi, Brend Zmijanej, "
$deliveryPara = $d.Paragraphs(26).Range
$fr = $d.Range($deliveryPara.Start, $deliveryPara.End)
$fr.Find.Execute(" -> Dea Hasanaj, Brend ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, " -> Danja Korreshi, Brend ", 2) | Out-Null

$deliveryPara = $d.Paragraphs(26).Range
$fr = $d.Range($deliveryPara.Start, $deliveryPara.End)
$fr.Find.Execute("Zmijanej", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "Zmijanej, ", 2) | Out-Null

# --- Restaurant html, css, js -> ... bullet -------------------------
# "Danja Korreshi, Aurel Kulemani, Olga Kolaj" => "Dea Hasanaj, Aurel Kulemani, Olga Kolaj"
$restaurantPara = $d.Paragraphs(27).Range
$fr = $d.Range($restaurantPara.Start, $restaurantPara.End)
$fr.Find.Execute(" -> Danja Korreshi, Aurel ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, " -> Dea Hasanaj, Aurel ", 2) | Out-Null

# --- Remove the now-empty spacer paragraph (indented 360 twips) that
#     used to follow the "Restaurant" bullet -------------------------
$spacer = $d.Paragraphs(28)
$spacer.Range.Delete()
